# Exceltabelle erweitert; Automatisch Endzeitdatum fuellen falls NULL war
#
# The existing last row (47) represented an open/unclosed attendance entry
# for "Stephan Fuchs" (Ankunft == previous Verlassen). A later scan for the
# same student closes/overwrites that open entry with the new timestamp and
# a fresh (reset) duration, and seven brand-new attendance rows are appended
# for the subsequent scan events recorded on 07.06.2024 and 08.06.2024.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 47: the new scan re-stamps both the "Ankunft"
#     and "Verlassen" timestamps of this entry to the new scan time,
#     resetting the duration for the new open session.
$ws.Range("D47").Value = "07.06.2024 16:16"
$ws.Range("E47").Value = "07.06.2024 16:16"
$ws.Range("F47").Value = 0

# --- Append new attendance rows 48-54 ---

# Row 48: Eli Enders, 4a - new arrival scan (still open, duration 0)
$ws.Range("A48").Value = "Eli"
$ws.Range("B48").Value = "Enders"
$ws.Range("C48").Value = "4a"
$ws.Range("D48").Value = "07.06.2024 16:16"
$ws.Range("E48").Value = "07.06.2024 16:16"
$ws.Range("F48").Value = 0

# Row 49: Max Schmitz, 4a - new arrival scan (still open, duration 0)
$ws.Range("A49").Value = "Max"
$ws.Range("B49").Value = "Schmitz"
$ws.Range("C49").Value = "4a"
$ws.Range("D49").Value = "07.06.2024 16:16"
$ws.Range("E49").Value = "07.06.2024 16:16"
$ws.Range("F49").Value = 0

# Row 50: Detlef Soost, 1a - new arrival scan (still open, duration 0)
$ws.Range("A50").Value = "Detlef"
$ws.Range("B50").Value = "Soost"
$ws.Range("C50").Value = "1a"
$ws.Range("D50").Value = "07.06.2024 16:16"
$ws.Range("E50").Value = "07.06.2024 16:16"
$ws.Range("F50").Value = 0

# Row 51: Stephan Fuchs, 3c - arrived and left on 08.06.2024
$ws.Range("A51").Value = "Stephan"
$ws.Range("B51").Value = "Fuchs"
$ws.Range("C51").Value = "3c"
$ws.Range("D51").Value = "08.06.2024 09:02"
$ws.Range("E51").Value = "08.06.2024 09:03"
$ws.Range("F51").Value = 1

# Row 52: Eli Enders, 4a - arrived and left on 08.06.2024
$ws.Range("A52").Value = "Eli"
$ws.Range("B52").Value = "Enders"
$ws.Range("C52").Value = "4a"
$ws.Range("D52").Value = "08.06.2024 09:02"
$ws.Range("E52").Value = "08.06.2024 09:06"
$ws.Range("F52").Value = 4

# Row 53: Detlef Soost, 1a - arrived and left on 08.06.2024
$ws.Range("A53").Value = "Detlef"
$ws.Range("B53").Value = "Soost"
$ws.Range("C53").Value = "1a"
$ws.Range("D53").Value = "08.06.2024 09:11"
$ws.Range("E53").Value = "08.06.2024 18:00"
$ws.Range("F53").Value = 529

# Row 54: Max Schmitz, 4a - arrived and left on 08.06.2024
$ws.Range("A54").Value = "Max"
$ws.Range("B54").Value = "Schmitz"
$ws.Range("C54").Value = "4a"
$ws.Range("D54").Value = "08.06.2024 09:11"
$ws.Range("E54").Value = "08.06.2024 18:00"
$ws.Range("F54").Value = 529
